# Re-remove "culture_collection" column from the MIxS template
# (INSDC2017 based re-check). The whole column AB (header
# "culture_collection", with its cell comment) is deleted, so every
# column to its right shifts one place to the left, and the comment
# that used to sit on the very last column (BM / wind_speed) goes away.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$startCol = 28   # AB
$endCol   = 65   # BM
$lastCol  = $endCol - 1   # BL, last column after the shift

# 1) Remember the comment text for every header cell from AB15 to BM15
#    (in column order) before anything else is touched. Legacy cell
#    comments in this workbook are not repositioned by cell/column
#    operations, so we will reassign their text ourselves afterwards.
$commentTexts = @()
for ($col = $startCol; $col -le $endCol; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $commentTexts += $cell.Comment.Text()
}

# 2) Shift the cell values/styles of AC15:BM15 one column to the left
#    (onto AB15:BL15) using a plain copy, which keeps column
#    definitions (widths) untouched, then blank out the now-duplicated
#    trailing cell (BM15).
$srcRange = $ws.Range("AC15:BM15")
$dstCell  = $ws.Range("AB15")
$srcRange.Copy($dstCell)
$ws.Range("BM15").ClearContents()

# 3) Stash the now-correct 64-column header row in a scratch row far
#    below the data so it survives a rebuild of row 15.
$ws.Range("A15:BL15").Copy($ws.Range("A1000"))

# 4) Rebuild row 15 from scratch (delete + insert is a row-axis
#    operation only, so it leaves column width metadata alone) to drop
#    the leftover, now-empty trailing cell and shrink the row's cell
#    range back down to 64 columns.
$ws.Rows(15).Delete()
$ws.Rows(15).Insert()

# 5) Restore the 64-column header row from the scratch copy, then
#    remove the scratch row again.
$ws.Range("A1000:BL1000").Copy($ws.Range("A15"))
$ws.Rows(1000).Delete()

# 6) Re-assign the comment text for AB15 .. BL15 using the values
#    captured in step 1, shifted left by one column to mirror what
#    happened to the cell contents.
for ($col = $startCol; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    [void]$cell.Comment.Text($commentTexts[$col - $startCol + 1])
}

# 7) The comment that used to belong to the last column (BM, the
#    "wind speed" one) is now orphaned with no backing cell content;
#    remove it entirely.
$orphan = $ws.Cells.Item($row, $endCol)
[void]$orphan.Comment.Delete()
